$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Replace the "1+1 gratis" promo formula with the new scalable heuristic strategy
$ws.Range("C6").Value = "if len(prices) >= 2: total_price -= sorted(prices)[1] * 0.3"
